$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("E3").Value = 16.314
$ws.Range("D12").Value = -7.328999999999999
$ws.Range("E14").Value = 16.876
$ws.Range("E26").Value = 16.946
$ws.Range("D27").Value = -8.684999999999999
$ws.Range("E31").Value = 17.064
$ws.Range("D32").Value = -8.175999999999998
$ws.Range("E35").Value = 16.63
$ws.Range("D36").Value = -7.783999999999999
$ws.Range("E37").Value = 16.737
$ws.Range("D38").Value = -7.662000000000001
$ws.Range("E45").Value = 16.96
$ws.Range("D46").Value = -8.23
$ws.Range("E52").Value = 16.768
$ws.Range("D54").Value = -8.42
$ws.Range("D55").Value = -8.129000000000001
$ws.Range("D56").Value = -8.34
$ws.Range("E57").Value = 16.643
$ws.Range("D67").Value = -7.467999999999999
$ws.Range("D69").Value = -7.555000000000001
$ws.Range("D72").Value = -7.595000000000001
$ws.Range("E81").Value = 16.858
$ws.Range("D83").Value = -7.966999999999999
$ws.Range("E83").Value = 16.661
$ws.Range("D86").Value = -8.263
$ws.Range("D91").Value = -7.417
$ws.Range("D93").Value = -7.855
$ws.Range("D99").Value = -8.019
$ws.Range("E100").Value = 16.684
$ws.Range("E102").Value = 16.669
